$d = $word.ActiveDocument
$d.Content.Find.MatchWildcards = $false

# 1. Title - no-op touch so Word re-serializes the paragraph as a single run.
$t1 = "Data Science 2: Statistics for Data Science"
$null = $d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2)

# 2. Subtitle year: 1871 - 2000 -> 1871 - 2020 (real text edit).
$null = $d.Content.Find.Execute("1871 - 2000", $true, $false, $false, $false, $false, $true, 1, $false, "1871 - 2020", 2)

# 3. "Submitted by Group 8" - no-op touch.
$t3 = "Submitted by Group 8"
$null = $d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2)

# 4. "Johanna Yu" - no-op touch.
$t4 = "Johanna Yu"
$null = $d.Content.Find.Execute($t4, $true, $false, $false, $false, $false, $true, 1, $false, $t4, 2)

# 5. "Dr. Amanda " - no-op touch.
$t5 = "Dr. Amanda "
$null = $d.Content.Find.Execute($t5, $true, $false, $false, $false, $false, $true, 1, $false, $t5, 2)

# 6. "December 2021" - no-op touch.
$t6 = "December 2021"
$null = $d.Content.Find.Execute($t6, $true, $false, $false, $false, $false, $true, 1, $false, $t6, 2)

# 7. "... confirm and/or exploit? " - no-op touch via wildcard so the embedded
#    Unicode right-single-quote (U+2019) is preserved without re-typing it.
$d.Content.Find.MatchWildcards = $true
$w7 = " being a correlation in the data that you*re looking to confirm and/or exploit\? "
$null = $d.Content.Find.Execute($w7, $true, $false, $true, $false, $false, $true, 1, $false, "^&", 2)
$d.Content.Find.MatchWildcards = $false

# 8. "Because we wanted..." paragraph - two real text edits.
$null = $d.Content.Find.Execute("we should minimize", $true, $false, $false, $false, $false, $true, 1, $false, "we wanted to minimize", 2)
$null = $d.Content.Find.Execute("the most purest data available.", $true, $false, $false, $false, $false, $true, 1, $false, "the purest unaltered data available.", 2)

# 9. Move the _GoBack bookmark into the middle of "player" in "the same player playing".
$found9 = $d.Content
$null = $found9.Find.Execute("same player playing")
$splitPos = $found9.Start + 7   # right after "same pl"
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 10. "diagrams to make your points. If you chose to do a predictive model, " - no-op touch.
$t10 = " diagrams to make your points. If you chose to do a predictive model, "
$null = $d.Content.Find.Execute($t10, $true, $false, $false, $false, $false, $true, 1, $false, $t10, 2)

# 11. "Did you prove/disprove..." - no-op touch to merge runs (also drops the
#     now-stale gramStart/gramEnd wrap around "your").
$t11 = "Did you prove/disprove your hypothesis or create a useful model? What did you learn about your data set?"
$null = $d.Content.Find.Execute($t11, $true, $false, $false, $false, $false, $true, 1, $false, $t11, 2)
